# Auto-generated Excel COM-interop script applying numeric updates
# to the Sargatanas_Profits leve-profit calculation sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 443.5
$ws.Range("I33").Value = 337.53488
$ws.Range("K33").Value = 337.53488
$ws.Range("M33").Value = -108.53488
$ws.Range("H70").Value = 34315624
$ws.Range("J70").Value = 30305522
$ws.Range("L70").Value = 90916566
$ws.Range("N70").Value = -90917106
$ws.Range("H73").Value = 34315624
$ws.Range("J73").Value = 30305522
$ws.Range("L73").Value = 90916566
$ws.Range("N73").Value = -90918438
$ws.Range("H107").Value = 26391764
$ws.Range("I107").Value = 11366836
$ws.Range("J107").Value = 50002370
$ws.Range("K107").Value = 11366836
$ws.Range("L107").Value = 50002370
$ws.Range("M107").Value = -11364916
$ws.Range("N107").Value = -50006210
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 12248.75
$ws.Range("I39").Value = 13666
$ws.Range("J39").Value = 7997
$ws.Range("K39").Value = 13666
$ws.Range("L39").Value = 7997
$ws.Range("M39").Value = -13146
$ws.Range("N39").Value = -9037
$ws.Range("H51").Value = 101666.664
$ws.Range("J51").Value = 110000
$ws.Range("L51").Value = 110000
$ws.Range("N51").Value = -111512
$ws.Range("H58").Value = 109999.664
$ws.Range("J58").Value = 109999.664
$ws.Range("L58").Value = 109999.664
$ws.Range("N58").Value = -110859.664
$ws.Range("H62").Value = 43000
$ws.Range("J62").Value = 43000
$ws.Range("L62").Value = 43000
$ws.Range("N62").Value = -44248
$ws.Range("H64").Value = 43000
$ws.Range("J64").Value = 43000
$ws.Range("L64").Value = 43000
$ws.Range("N64").Value = -43496
$ws.Range("H65").Value = 43000
$ws.Range("J65").Value = 43000
$ws.Range("L65").Value = 129000
$ws.Range("N65").Value = -135240
$ws.Range("H67").Value = 43000
$ws.Range("J67").Value = 43000
$ws.Range("L67").Value = 43000
$ws.Range("N67").Value = -44716
$ws.Range("H76").Value = 65000
$ws.Range("J76").Value = 65000
$ws.Range("L76").Value = 65000
$ws.Range("N76").Value = -65676
$ws.Range("H79").Value = 65000
$ws.Range("J79").Value = 65000
$ws.Range("L79").Value = 65000
$ws.Range("N79").Value = -67340

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 37468.75
$ws.Range("J18").Value = 37468.75
$ws.Range("L18").Value = 37468.75
$ws.Range("N18").Value = -37928.75
$ws.Range("H31").Value = 9650.083000000001
$ws.Range("I31").Value = 3954.1428
$ws.Range("J31").Value = 17624.4
$ws.Range("K31").Value = 3954.1428
$ws.Range("L31").Value = 17624.4
$ws.Range("M31").Value = -3659.1428
$ws.Range("N31").Value = -18214.4
$ws.Range("H34").Value = 9650.083000000001
$ws.Range("I34").Value = 3954.1428
$ws.Range("J34").Value = 17624.4
$ws.Range("K34").Value = 3954.1428
$ws.Range("L34").Value = 17624.4
$ws.Range("M34").Value = -3752.1428
$ws.Range("N34").Value = -18028.4
$ws.Range("H36").Value = 28999
$ws.Range("J36").Value = 28999
$ws.Range("L36").Value = 28999
$ws.Range("N36").Value = -29775
$ws.Range("H40").Value = 28999
$ws.Range("J40").Value = 28999
$ws.Range("L40").Value = 28999
$ws.Range("N40").Value = -29319
$ws.Range("H42").Value = 19994
$ws.Range("J42").Value = 19994
$ws.Range("L42").Value = 19994
$ws.Range("N42").Value = -21180
$ws.Range("H53").Value = 290490
$ws.Range("J53").Value = 290490
$ws.Range("L53").Value = 290490
$ws.Range("N53").Value = -291704
$ws.Range("H74").Value = 125047420
$ws.Range("J74").Value = 54199.57
$ws.Range("L74").Value = 54199.57
$ws.Range("N74").Value = -55947.57
$ws.Range("H77").Value = 125047420
$ws.Range("J77").Value = 54199.57
$ws.Range("L77").Value = 162598.71
$ws.Range("N77").Value = -171334.71
$ws.Range("H107").Value = 1371.9445
$ws.Range("I107").Value = 621.7143
$ws.Range("K107").Value = 621.7143
$ws.Range("M107").Value = 1298.2857
$ws.Range("H134").Value = 3569.88
$ws.Range("I134").Value = 1692.738
$ws.Range("K134").Value = 5078.214
$ws.Range("M134").Value = -2543.214

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2622.739
$ws.Range("J68").Value = 2734.9443
$ws.Range("L68").Value = 8204.832900000001
$ws.Range("N68").Value = -9826.832900000001
$ws.Range("H71").Value = 2622.739
$ws.Range("J71").Value = 2734.9443
$ws.Range("L71").Value = 24614.4987
$ws.Range("N71").Value = -32726.4987
$ws.Range("H131").Value = 2048.261
$ws.Range("I131").Value = 1176.1111
$ws.Range("K131").Value = 3528.3333
$ws.Range("M131").Value = 1511.6667
$ws.Range("H137").Value = 89986
$ws.Range("I137").Value = 78051.766
$ws.Range("K137").Value = 234155.298
$ws.Range("M137").Value = -229055.298

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 15916.667
$ws.Range("J36").Value = 19750
$ws.Range("L36").Value = 19750
$ws.Range("N36").Value = -20720
$ws.Range("H80").Value = 2829.818
$ws.Range("I80").Value = 1963
$ws.Range("K80").Value = 1963
$ws.Range("M80").Value = -965
$ws.Range("H83").Value = 2829.818
$ws.Range("I83").Value = 1963
$ws.Range("K83").Value = 9815
$ws.Range("M83").Value = -4823
$ws.Range("H97").Value = 1864.8823
$ws.Range("I97").Value = 1913.3125
$ws.Range("K97").Value = 1913.3125
$ws.Range("M97").Value = -1417.3125
$ws.Range("H122").Value = 33370502
$ws.Range("I122").Value = 58885550
$ws.Range("K122").Value = 176656650
$ws.Range("M122").Value = -176654200
$ws.Range("H132").Value = 5307.9473
$ws.Range("I132").Value = 3340.7
$ws.Range("K132").Value = 10022.1
$ws.Range("M132").Value = -7492.099999999999
$ws.Range("H141").Value = 59999
$ws.Range("J141").Value = 59999
$ws.Range("L141").Value = 59999
$ws.Range("N141").Value = -70359

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1402.7778
$ws.Range("I22").Value = 1076.1538
$ws.Range("J22").Value = 2252
$ws.Range("K22").Value = 1076.1538
$ws.Range("L22").Value = 2252
$ws.Range("M22").Value = -781.1538
$ws.Range("N22").Value = -2842
$ws.Range("H27").Value = 1402.7778
$ws.Range("I27").Value = 1076.1538
$ws.Range("J27").Value = 2252
$ws.Range("K27").Value = 1076.1538
$ws.Range("L27").Value = 2252
$ws.Range("M27").Value = -969.1538
$ws.Range("N27").Value = -2466
$ws.Range("H100").Value = 7000
$ws.Range("I100").Value = 7000
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 7000
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -6459
$ws.Range("N100").ClearContents()
$ws.Range("H122").Value = 6414.711
$ws.Range("I122").Value = 6461.241
$ws.Range("J122").Value = 6330.375
$ws.Range("K122").Value = 19383.723
$ws.Range("L122").Value = 18991.125
$ws.Range("M122").Value = -16933.723
$ws.Range("N122").Value = -23891.125
$ws.Range("H132").Value = 6196.154
$ws.Range("I132").Value = 3567.6956
$ws.Range("K132").Value = 10703.0868
$ws.Range("M132").Value = -8173.086800000001
$ws.Range("H136").Value = 12543.634
$ws.Range("I136").Value = 10830.929
$ws.Range("J136").Value = 14042.25
$ws.Range("K136").Value = 32492.787
$ws.Range("L136").Value = 42126.75
$ws.Range("M136").Value = -29942.787
$ws.Range("N136").Value = -47226.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 13827.579
$ws.Range("I113").Value = 21215.916
$ws.Range("K113").Value = 63647.74800000001
$ws.Range("M113").Value = -61477.74800000001
$ws.Range("H122").Value = 7528227
$ws.Range("I122").Value = 10503981
$ws.Range("K122").Value = 31511943
$ws.Range("M122").Value = -31509493
$ws.Range("H132").Value = 12724.291
$ws.Range("I132").Value = 9436.23
$ws.Range("K132").Value = 28308.69
$ws.Range("M132").Value = -25778.69

